$d = $word.ActiveDocument

# --- Step 1: Replace " Started writing code for " run text with " " ---
$rFind = $d.Content.Duplicate
$rFind.Find.Execute(" Started writing code for ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rFind.Delete()
$rFind.InsertAfter(" ")

# --- Step 2: Remove the _GoBack bookmark (will re-add at new end) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3: Insert "Wrote basic host-side code for " run right after the " " run we just set, before "key issuing system..." ---
$rFind.Collapse(0)
$rFind.InsertAfter("Wrote basic host-side code for ")

# --- Step 4: append the rest of paragraph 0 runs (after "key issuing system...APDU.") at paragraph end ---
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter(" Created new compile/install scripts")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)

# --- Step 5: append new paragraphs ---
# paragraph 1
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Attempted to compile. Various errors relating to uses of, and indirect casts to, ")
$r.Collapse(0)
$r.InsertAfter("int")
$r.Collapse(0)
$r.InsertAfter(" instead of short. Fixed those.")
$r.Collapse(0)
$r.InsertAfter(" Successfully compiles and converts. ")
$r.Collapse(0)
$r.InsertAfter("GPShell")
$r.Collapse(0)
$r.InsertAfter(" install command returns 6A80")
$r.Collapse(0)
$r.InsertAfter(" (wrong data / incorrect values in command data)")
$r.Collapse(0)
$r.InsertAfter(" Will fix tomorrow.")
$r.Collapse(0)
# paragraph 2
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("30/11 Fixed. Problem was likely due to install method not calling register.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Kp.genKeyPair")
$r.Collapse(0)
$r.InsertAfter("() not working, throws ")
$r.Collapse(0)
$r.InsertAfter("(6F, 00)")
$r.Collapse(0)
$r.InsertAfter(". Turns out it was because ")
$r.Collapse(0)
$r.InsertAfter("KeyPair")
$r.Collapse(0)
$r.InsertAfter(" constructor can’t take arbitrary key lengths, only the constants in ")
$r.Collapse(0)
$r.InsertAfter("KeyBuilder")
$r.Collapse(0)
$r.InsertAfter(". In JC 2.2.2 they only go to 192b, but protocol requires 256b. Instead, ")
$r.Collapse(0)
$r.InsertAfter("have to")
$r.Collapse(0)
$r.InsertAfter(" separately initialise public and private keys using NIST EC parameters and use the other ")
$r.Collapse(0)
$r.InsertAfter("KeyPair")
$r.Collapse(0)
$r.InsertAfter(" constructor.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("After adding it in, found it didn’t accept the compressed G, so had to enter the full uncompressed version. Fixed this, ")
$r.Collapse(0)
$r.InsertAfter("now generated keypair without error.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Had trouble with ")
$r.Collapse(0)
$r.InsertAfter("getW")
$r.Collapse(0)
$r.InsertAfter("(")
$r.Collapse(0)
$r.InsertAfter(") command to obtain public key. After a while figured out 64B buffer not big enough.")
$r.Collapse(0)
$r.InsertAfter(" Need 65B to accommodate extra 0x04 at the beginning.")
$r.Collapse(0)
# paragraph 3
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Am now reasonably confident key issuing process works. Will move on to debugging authentication process.")
$r.Collapse(0)
# paragraph 4
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("3/12 Fixed usage of hash function and other errors on the authentication code. Removed unnecessary array conversion functions and cleaned up code.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Implemented various uncompleted functions on host side including ")
$r.Collapse(0)
$r.InsertAfter("ec_dh")
$r.Collapse(0)
$r.InsertAfter(", ")
$r.Collapse(0)
$r.InsertAfter("kdf")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
# paragraph 5
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Realised that ")
$r.Collapse(0)
$r.InsertAfter("assumption behind signature was wrong. Should be calculated by the card and sent to host. Adapted issuing code to account for this.")
$r.Collapse(0)

# --- Step 6: re-add the _GoBack bookmark at the very end, collapsed (zero width) ---
$r.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $r)
$bmFinal = $d.Bookmarks("_GoBack")
$bmFinal.Range.Text = ""

Write-Output "Done"